# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The sheet holds a small metadata table (4 header rows x 11 columns) that
# documents each data column (measure/dimension, role, datatype/URI). The
# dimensions for "municipio" (C), "aragon"/comunidad (E), "situacion-profesional"
# (F) and "sexo" (G) are re-curated:
#   - "municipio" (col C) becomes a plain refArea dimension with a URI-Municipio
#   - "aragon" (col E) becomes a plain refArea dimension with a URI-Comunidad
#   - "situacion-profesional" (col F) and "sexo" (col G) stop being
#     skos:Concept dimensions and become plain int measures instead
# The now-obsolete mapping-file row (row 5) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (municipio-nombre): was an iaest-measure -> now a refArea dimension
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# --- Column E (aragon / comunidad): was an iaest-dimension with skos:Concept -> now a refArea dimension
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Comunidad"

# --- Column F (situacion-profesional): was a skos:Concept dimension -> now a plain measure
$ws.Range("F2").Value = "iaest-measure:situacion-profesional"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"

# --- Column G (sexo): was a skos:Concept dimension -> now a plain measure
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"

# Row 5 (the old mapping-*.xlsx file references for aragon/situacion-profesional/sexo)
# is no longer needed, since those columns are no longer skos:Concept dimensions.
$ws.Range("E5:G5").Clear()
